# Scheduled-runner refresh of market-board pricing for the leve profit sheets.
# For each affected leve row, re-derive the currentAveragePrice* (H:J),
# LevePrice* (K:L) and LeveProfit* (M:N) columns from the latest pull.
$wb = $excel.ActiveWorkbook

# ALC!32 - Automata for the People / Crab Oil
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 609.06665
$ws.Range("I32").Value = 733.3333
$ws.Range("J32").Value = 578
$ws.Range("K32").Value = 733.3333
$ws.Range("L32").Value = 578
$ws.Range("M32").Value = -407.3333
$ws.Range("N32").Value = -1230

# ALC!98 - The Dotted Line / Enchanted Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1915.7234
$ws.Range("I98").Value = 1305.5581
$ws.Range("J98").Value = 8475
$ws.Range("K98").Value = 1305.5581
$ws.Range("L98").Value = 8475
$ws.Range("M98").Value = 192.4419
$ws.Range("N98").Value = -11471

# ALC!99 - Rumor Has It / Commanding Craftsman's Tea
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1334.5
$ws.Range("I99").Value = 913.5
$ws.Range("J99").Value = 2316.8333
$ws.Range("K99").Value = 2740.5
$ws.Range("L99").Value = 6950.499899999999
$ws.Range("M99").Value = -1242.5
$ws.Range("N99").Value = -9946.499899999999

# ALC!122 - Wishful Inking / Enchanted High Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1915.7234
$ws.Range("I122").Value = 1305.5581
$ws.Range("J122").Value = 8475
$ws.Range("K122").Value = 3916.6743
$ws.Range("L122").Value = 25425
$ws.Range("M122").Value = -1466.6743
$ws.Range("N122").Value = -30325

# ALC!125 - Body over Mind / Grade 5 Dexterity Alkahest
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 20802
$ws.Range("I125").Value = 40737.332
$ws.Range("J125").Value = 866.6667
$ws.Range("K125").Value = 366635.988
$ws.Range("L125").Value = 7800.0003
$ws.Range("M125").Value = -364175.988
$ws.Range("N125").Value = -12720.0003

# ALC!127 - Liquid Competence / Competent Craftsman's Draught
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1384.909
$ws.Range("I127").Value = 1316.6666
$ws.Range("J127").Value = 1410.5
$ws.Range("K127").Value = 3949.9998
$ws.Range("L127").Value = 4231.5
$ws.Range("M127").Value = 1010.0002
$ws.Range("N127").Value = -14151.5

# ALC!137 - Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3226.1965
$ws.Range("I137").Value = 1063.303
$ws.Range("J137").Value = 6329.478
$ws.Range("K137").Value = 3189.909000000001
$ws.Range("L137").Value = 18988.434
$ws.Range("M137").Value = -639.9090000000006
$ws.Range("N137").Value = -24088.434

# ARM!15 - All Ovo That / Iron Skillet
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 14000
$ws.Range("J15").Value = 14000
$ws.Range("L15").Value = 14000
$ws.Range("N15").Value = -14700

# ARM!21 - Fashion Weak / Iron Cuirass
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1859.1428
$ws.Range("I21").Value = 2319
$ws.Range("J21").Value = 709.5
$ws.Range("K21").Value = 2319
$ws.Range("L21").Value = 709.5
$ws.Range("M21").Value = -1945
$ws.Range("N21").Value = -1457.5

# ARM!61 - Dealing with the Tough Stuff / Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 272056.88
$ws.Range("I61").Value = 1789.8572
$ws.Range("J61").Value = 1112887.5
$ws.Range("K61").Value = 1789.8572
$ws.Range("L61").Value = 1112887.5
$ws.Range("M61").Value = -1577.8572
$ws.Range("N61").Value = -1113311.5

# ARM!132 - Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4659.8687
$ws.Range("I132").Value = 3408.9768
$ws.Range("K132").Value = 10226.9304
$ws.Range("M132").Value = -7696.930399999999

# ARM!136 - Metal with Mettle / Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 272056.88
$ws.Range("I136").Value = 1789.8572
$ws.Range("J136").Value = 1112887.5
$ws.Range("K136").Value = 5369.571599999999
$ws.Range("L136").Value = 3338662.5
$ws.Range("M136").Value = -2819.571599999999
$ws.Range("N136").Value = -3343762.5

# BSM!134 - Ruthenium Supremium / Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1214.409
$ws.Range("I134").Value = 681.8823
$ws.Range("J134").Value = 3025
$ws.Range("K134").Value = 2045.6469
$ws.Range("L134").Value = 9075
$ws.Range("M134").Value = 489.3531
$ws.Range("N134").Value = -14145

# CRP!22 - Driving Up the Wall / Elm Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 832.25
$ws.Range("I22").Value = 132.28572
$ws.Range("J22").Value = 1812.2
$ws.Range("K22").Value = 132.28572
$ws.Range("L22").Value = 1812.2
$ws.Range("M22").Value = 217.71428
$ws.Range("N22").Value = -2512.2

# CRP!51 - Greenstone for Greenhorns / Jade Crook
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 26198.572
$ws.Range("I51").Value = 30000
$ws.Range("J51").Value = 25565
$ws.Range("K51").Value = 30000
$ws.Range("L51").Value = 25565
$ws.Range("M51").Value = -29264
$ws.Range("N51").Value = -27037

# CRP!58 - You Do the Heavy Lifting / Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 27778834
$ws.Range("I58").Value = 50000650
$ws.Range("J58").Value = 1566.75
$ws.Range("K58").Value = 50000650
$ws.Range("L58").Value = 1566.75
$ws.Range("M58").Value = -50000447
$ws.Range("N58").Value = -1972.75

# CRP!61 - Incant Now, Think Later / Jade Crook
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 26198.572
$ws.Range("I61").Value = 30000
$ws.Range("J61").Value = 25565
$ws.Range("K61").Value = 30000
$ws.Range("L61").Value = 25565
$ws.Range("M61").Value = -29652
$ws.Range("N61").Value = -26261

# CRP!64 - Almost as Fun as Slingshotting Birds / Cedar Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 27769.23
$ws.Range("J64").Value = 27769.23
$ws.Range("L64").Value = 27769.23
$ws.Range("N64").Value = -28265.23

# CRP!67 - Living Bow to Mouth (L) / Cedar Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 27769.23
$ws.Range("J67").Value = 27769.23
$ws.Range("L67").Value = 27769.23
$ws.Range("N67").Value = -29485.23

# CRP!105 - Zelkova, My Love / Zelkova Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 10403.667
$ws.Range("I105").Value = 10227.75
$ws.Range("J105").Value = 10755.5
$ws.Range("K105").Value = 10227.75
$ws.Range("L105").Value = 10755.5
$ws.Range("M105").Value = -8480.75
$ws.Range("N105").Value = -14249.5

# CRP!118 - A Miss and a Hit / Sandteak Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# CRP!136 - Turali Quality / Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 27778834
$ws.Range("I136").Value = 50000650
$ws.Range("J136").Value = 1566.75
$ws.Range("K136").Value = 150001950
$ws.Range("L136").Value = 4700.25
$ws.Range("M136").Value = -149999400
$ws.Range("N136").Value = -9800.25

# CUL!5 - What a Sap / Maple Syrup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 703.64703
$ws.Range("I5").Value = 496.85
$ws.Range("J5").Value = 999.0714
$ws.Range("K5").Value = 1490.55
$ws.Range("L5").Value = 2997.2142
$ws.Range("M5").Value = -1378.55
$ws.Range("N5").Value = -3221.2142

# CUL!107 - Slippery Service / Frantoio Oil
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1043.8422
$ws.Range("I107").Value = 1064.7778
$ws.Range("J107").Value = 1025
$ws.Range("K107").Value = 3194.3334
$ws.Range("L107").Value = 3075
$ws.Range("M107").Value = -1274.3334
$ws.Range("N107").Value = -6915

# CUL!109 - Cure for What Ails / Purple Carrot Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 921.8570999999999
$ws.Range("I109").Value = 492.16666
$ws.Range("J109").Value = 3500
$ws.Range("K109").Value = 1476.49998
$ws.Range("L109").Value = 10500
$ws.Range("M109").Value = -436.4999800000001
$ws.Range("N109").Value = -12580

# CUL!114 - One Last Meal / Mushroom Saute
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1341.1333
$ws.Range("I114").Value = 508
$ws.Range("J114").Value = 2070.125
$ws.Range("K114").Value = 1524
$ws.Range("L114").Value = 6210.375
$ws.Range("M114").Value = 1730
$ws.Range("N114").Value = -12718.375

# CUL!117 - A Good Omen / Peppered Popotoes
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 440
$ws.Range("I117").Value = 413.33334
$ws.Range("J117").Value = 466.66666
$ws.Range("K117").Value = 1240.00002
$ws.Range("L117").Value = 1399.99998
$ws.Range("M117").Value = 2201.99998
$ws.Range("N117").Value = -8283.999980000001

# CUL!121 - A Cookie for Your Troubles / Coffee Biscuit
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 8092
$ws.Range("I121").Value = 593.625
$ws.Range("J121").Value = 9670.605
$ws.Range("K121").Value = 1780.875
$ws.Range("L121").Value = 29011.815
$ws.Range("M121").Value = -470.875
$ws.Range("N121").Value = -31631.815

# CUL!122 - Salt of the North / Northern Sea Salt
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 809.3095
$ws.Range("I122").Value = 414.9524
$ws.Range("J122").Value = 1203.6666
$ws.Range("K122").Value = 3734.5716
$ws.Range("L122").Value = 10832.9994
$ws.Range("M122").Value = -1284.5716
$ws.Range("N122").Value = -15732.9994

# CUL!131 - The Mountain Steeped / Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 873.5208
$ws.Range("I131").Value = 481.33334
$ws.Range("J131").Value = 964.02563
$ws.Range("K131").Value = 1444.00002
$ws.Range("L131").Value = 2892.07689
$ws.Range("M131").Value = 3595.99998
$ws.Range("N131").Value = -12972.07689

# CUL!135 - Not-so-secret Ingredient / Royal Maple Syrup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 703.64703
$ws.Range("I135").Value = 496.85
$ws.Range("J135").Value = 999.0714
$ws.Range("K135").Value = 4471.650000000001
$ws.Range("L135").Value = 8991.642600000001
$ws.Range("M135").Value = -1936.650000000001
$ws.Range("N135").Value = -14061.6426

# CUL!140 - Sweet, Sweet Bean Juice / Mesquite Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 29601.225
$ws.Range("I140").Value = 36201.562
$ws.Range("J140").Value = 3199.875
$ws.Range("K140").Value = 108604.686
$ws.Range("L140").Value = 9599.625
$ws.Range("M140").Value = -103424.686
$ws.Range("N140").Value = -19959.625

# GSM!122 - Awarding Academic Excellence / Ametrine
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2633.3928
$ws.Range("I122").Value = 3704.75
$ws.Range("J122").Value = 1829.875
$ws.Range("K122").Value = 11114.25
$ws.Range("L122").Value = 5489.625
$ws.Range("M122").Value = -8664.25
$ws.Range("N122").Value = -10389.625

# GSM!126 - Gold Rush Order / Phrygian Gold Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1554.7142
$ws.Range("I126").Value = 1371.4286
$ws.Range("J126").Value = 1921.2858
$ws.Range("K126").Value = 4114.2858
$ws.Range("L126").Value = 5763.857400000001
$ws.Range("M126").Value = -1644.2858
$ws.Range("N126").Value = -10703.8574

# GSM!132 - On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4683.878
$ws.Range("I132").Value = 5815.3706
$ws.Range("J132").Value = 2501.7144
$ws.Range("K132").Value = 17446.1118
$ws.Range("L132").Value = 7505.1432
$ws.Range("M132").Value = -14916.1118
$ws.Range("N132").Value = -12565.1432

# LTW!22 - Skin off Their Backs / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 361.2353
$ws.Range("I22").Value = 261.57144
$ws.Range("J22").Value = 431
$ws.Range("K22").Value = 261.57144
$ws.Range("L22").Value = 431
$ws.Range("M22").Value = 33.42856
$ws.Range("N22").Value = -1021

# LTW!27 - Fire and Hide / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 361.2353
$ws.Range("I27").Value = 261.57144
$ws.Range("J27").Value = 431
$ws.Range("K27").Value = 261.57144
$ws.Range("L27").Value = 431
$ws.Range("M27").Value = -154.57144
$ws.Range("N27").Value = -645

# LTW!30 - Packing a Punch / Goatskin Cesti
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 900
$ws.Range("I30").Value = 900
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 900
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -792
$ws.Range("N30").ClearContents()

# LTW!46 - Supply Side Logic / Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1438.7778
$ws.Range("I46").Value = 1029.8
$ws.Range("J46").Value = 1950
$ws.Range("K46").Value = 1029.8
$ws.Range("L46").Value = 1950
$ws.Range("M46").Value = -841.8
$ws.Range("N46").Value = -2326

# LTW!55 - It's Not a Job, It's a Calling / Peiste Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 222.95238
$ws.Range("I55").Value = 187.84616
$ws.Range("J55").Value = 280
$ws.Range("K55").Value = 187.84616
$ws.Range("L55").Value = 280
$ws.Range("M55").Value = -14.84616
$ws.Range("N55").Value = -626
